$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Narrow a handful of columns from width 8 to width 7 ---
# (columns G, K, P, V -> index 7, 11, 16, 22)
# Note: the raw OOXML "width" attribute is ColumnWidth + 0.83 in this font,
# so ColumnWidth = 6.17 serializes to a stored width of 7.
$ws.Columns.Item(7).ColumnWidth = 6.17
$ws.Columns.Item(11).ColumnWidth = 6.17
$ws.Columns.Item(16).ColumnWidth = 6.17
$ws.Columns.Item(22).ColumnWidth = 6.17

# --- Apply "custom accuracy" rounding to row 5 data values ---
$ws.Range("B5").Value = 14.91
$ws.Range("C5").Value = 11.45
$ws.Range("D5").Value = 0.45
$ws.Range("E5").Value = 32.04
$ws.Range("F5").Value = 26.86
$ws.Range("G5").Value = 11.38
$ws.Range("H5").Value = 42.78
$ws.Range("I5").Value = 17.78
$ws.Range("J5").Value = 8
$ws.Range("K5").Value = 11.82
$ws.Range("L5").Value = 13.38
$ws.Range("M5").Value = 13.92
$ws.Range("N5").Value = 3.83
$ws.Range("O5").Value = 11.5
$ws.Range("P5").Value = 16.82
$ws.Range("Q5").Value = 9.64
$ws.Range("R5").Value = 0.43
$ws.Range("S5").Value = 0.31
$ws.Range("T5").Value = 169.86
$ws.Range("U5").Value = 32.54
$ws.Range("V5").Value = 10.87
$ws.Range("W5").Value = 21.85
$ws.Range("X5").Value = 11.41
$ws.Range("Y5").Value = 1.52
$ws.Range("Z5").Value = 21.34
$ws.Range("AA5").Value = 9.5
$ws.Range("AB5").Value = 8.2
$ws.Range("AC5").Value = 9.85
$ws.Range("AD5").Value = 13.61
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 39.04
$ws.Range("AG5").Value = 6.65
$ws.Range("AH5").Value = 13.28

# --- Remove row 6 entirely (data trimmed from the sheet) ---
$ws.Rows.Item(6).Delete()
